# Update countries & provincias Spain
#
# The "Pais" sheet is a COVID-19 dashboard ranking table sorted descending
# by "Casos totales" (column B). The underlying stats for a handful of
# countries were refreshed; two of them (Bielorrusia / Corea del Sur and
# Sri Lanka / San Marino) overtook their neighbour in the ranking once the
# new totals were applied, so those two row-pairs swap places (country name
# + all of that row's stats move together) while every other affected row
# simply gets new numbers in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 - Rumania: refreshed stats (stays in place)
$ws.Range("B34").Value = 11339
$ws.Range("C34").Value = 303
$ws.Range("D34").Value = 3141
$ws.Range("E34").Value = 7567
$ws.Range("F34").Value = 227

# Row 35 - now Bielorrusia (was Corea del Sur): moved up, new stats
$ws.Range("A35").Value = "Bielorrusia"
$ws.Range("B35").Value = 11289
$ws.Range("C35").Value = 826
$ws.Range("D35").Value = 1740
$ws.Range("E35").Value = 9474
$ws.Range("F35").Value = 92
$ws.Range("G35").Value = 3
$ws.Range("H35").Value = 75

# Row 36 - now Corea del Sur (was Bielorrusia): moved down, keeps its
# previous (unchanged) stats
$ws.Range("A36").Value = "Corea del Sur"
$ws.Range("B36").Value = 10738
$ws.Range("C36").Value = 10
$ws.Range("D36").Value = 8764
$ws.Range("E36").Value = 1731
$ws.Range("F36").Value = 55
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = 243

# Row 52 - Finlandia: refreshed stats
$ws.Range("B52").Value = 4695
$ws.Range("C52").Value = 119
$ws.Range("E52").Value = 2005

# Row 55 - Marruecos: refreshed stats
$ws.Range("B55").Value = 4115
$ws.Range("C55").Value = 50
$ws.Range("D55").Value = 669
$ws.Range("E55").Value = 3285

# Row 68 - Uzbekistan: refreshed stats
$ws.Range("D68").Value = 826
$ws.Range("E68").Value = 1053

# Row 87 - Hong Kong: refreshed stats
$ws.Range("D87").Value = 787
$ws.Range("E87").Value = 247

# Row 105 - now Sri Lanka (was San Marino): moved up, new stats
$ws.Range("A105").Value = "Sri Lanka"
$ws.Range("B105").Value = 557
$ws.Range("C105").Value = 34
$ws.Range("D105").Value = 126
$ws.Range("E105").Value = 424
$ws.Range("F105").Value = 2
$ws.Range("H105").Value = 7

# Row 106 - now San Marino (was Sri Lanka): moved down, keeps its
# previous (unchanged) stats
$ws.Range("A106").Value = "San Marino"
$ws.Range("B106").Value = 538
$ws.Range("D106").Value = 64
$ws.Range("E106").Value = 433
$ws.Range("F106").Value = 4
$ws.Range("H106").Value = 41

# Row 117 - Kenia: refreshed stats
$ws.Range("B117").Value = 363
$ws.Range("C117").Value = 8
$ws.Range("E117").Value = 243
